$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30..114 down to 31..115
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record's data
$ws.Cells.Item(30, 1).Value = 4
$ws.Cells.Item(30, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value = "Los Lagos"
$ws.Cells.Item(30, 4).Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(44659)
$ws.Cells.Item(30, 5).Value = 10
$ws.Cells.Item(30, 6).Value = 100112022
$ws.Cells.Item(30, 7).Value = "Arveja Verde"
$ws.Cells.Item(30, 8).Value = "Perfection"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 70
$ws.Cells.Item(30, 11).Value = 33000
$ws.Cells.Item(30, 12).Value = 33000
$ws.Cells.Item(30, 13).Value = 33000
$ws.Cells.Item(30, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(30, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(30, 16).Value = 1320
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
